$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.179728507995605
$ws.Range("B1").Value = 1.567793726921082
$ws.Range("C1").Value = 1.47925591468811
$ws.Range("D1").Value = 1.433712244033813
$ws.Range("E1").Value = 1.314712882041931
